$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Creatures")

# --- Row 166: Aralez ---
$ws.Range("A166").Value = "Aralez"
$ws.Range("B166").Value = "1.0.0"
$ws.Range("C166").Value = 3
$ws.Range("D166").Value = "Celestial"
$ws.Range("E166").Value = "Arctic, Coastal, Desert, Extraplanar, Forest, Freshwater, Grassland, Hills, Jungle, Mountain, Swamp, Urban"
$ws.Range("G166").Value = "None"
$ws.Range("H166").Value = "None"
$ws.Range("I166").Value = "None"

# --- Row 167: Chamrosh ---
$ws.Range("A167").Value = "Chamrosh"
$ws.Range("B167").Value = "1.0.0"
$ws.Range("C167").Value = 7
$ws.Range("D167").Value = "Celestial"
$ws.Range("E167").Value = "Coastal, Desert, Extraplanar, Forest, Freshwater, Grassland, Hills, Jungle, Mountain, Swamp,"
$ws.Range("G167").Value = "None"
$ws.Range("H167").Value = "None"
$ws.Range("I167").Value = "None"

# --- Row 168: Flood Hound ---
$ws.Range("A168").Value = "Flood Hound"
$ws.Range("B168").Value = "1.0.0"
$ws.Range("C168").Value = 1
$ws.Range("D168").Value = "Elemental"
$ws.Range("E168").Value = "Arctic, Coastal, Desert, Extraplanar, Forest, Freshwater, Grassland, Hills, Jungle, Mountain, Swamp, Underwater, Urban"
$ws.Range("G168").Value = "None"
$ws.Range("H168").Value = "None"
$ws.Range("I168").Value = "None"

# --- Row 169: Miniature Blink Dog ---
$ws.Range("A169").Value = "Miniature Blink Dog"
$ws.Range("B169").Value = "1.0.0"
$ws.Range("C169").Value = 0.125
$ws.Range("D169").Value = "Fey"
$ws.Range("E169").Value = "Extraplanar, Forest, Grassland, Hills"
$ws.Range("G169").Value = "None"
$ws.Range("H169").Value = "None"
$ws.Range("I169").Value = "None"

# --- Column F (Public Source Doc), filled last, with hyperlinks ---
$ws.Range("F166").Value = "Magic Dogs"
$ws.Range("F167").Value = "Magic Dogs"
$ws.Range("F168").Value = "Magic Dogs"
$ws.Range("F169").Value = "Magic Dogs"

$ws.Hyperlinks.Add($ws.Range("F166"), "https://www.patreon.com/posts/magic-dogs-76221436")
$ws.Hyperlinks.Add($ws.Range("F167:F169"), "https://www.patreon.com/posts/magic-dogs-76221436", "", "", "Magic Dogs")

# Restore / normalize formatting for the whole new block so new cells reuse
# the same style indices as the rows directly above them (avoids creating
# duplicate style entries that a manual Excel edit wouldn't create).
$ws.Range("A164:I164").Copy()
$ws.Range("A166:I166").PasteSpecial(-4122)
$ws.Range("A164:I164").Copy()
$ws.Range("A167:I167").PasteSpecial(-4122)
$ws.Range("A164:I164").Copy()
$ws.Range("A168:I168").PasteSpecial(-4122)
$ws.Range("A164:I164").Copy()
$ws.Range("A169:I169").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply values, since PasteSpecial(formats) shouldn't touch them, but
# make sure the hyperlink text is exactly right for each cell.
$ws.Range("A166").Value = "Aralez"
$ws.Range("B166").Value = "1.0.0"
$ws.Range("C166").Value = 3
$ws.Range("D166").Value = "Celestial"
$ws.Range("E166").Value = "Arctic, Coastal, Desert, Extraplanar, Forest, Freshwater, Grassland, Hills, Jungle, Mountain, Swamp, Urban"
$ws.Range("F166").Value = "Magic Dogs"
$ws.Range("G166").Value = "None"
$ws.Range("H166").Value = "None"
$ws.Range("I166").Value = "None"

$ws.Range("A167").Value = "Chamrosh"
$ws.Range("B167").Value = "1.0.0"
$ws.Range("C167").Value = 7
$ws.Range("D167").Value = "Celestial"
$ws.Range("E167").Value = "Coastal, Desert, Extraplanar, Forest, Freshwater, Grassland, Hills, Jungle, Mountain, Swamp,"
$ws.Range("F167").Value = "Magic Dogs"
$ws.Range("G167").Value = "None"
$ws.Range("H167").Value = "None"
$ws.Range("I167").Value = "None"

$ws.Range("A168").Value = "Flood Hound"
$ws.Range("B168").Value = "1.0.0"
$ws.Range("C168").Value = 1
$ws.Range("D168").Value = "Elemental"
$ws.Range("E168").Value = "Arctic, Coastal, Desert, Extraplanar, Forest, Freshwater, Grassland, Hills, Jungle, Mountain, Swamp, Underwater, Urban"
$ws.Range("F168").Value = "Magic Dogs"
$ws.Range("G168").Value = "None"
$ws.Range("H168").Value = "None"
$ws.Range("I168").Value = "None"

$ws.Range("A169").Value = "Miniature Blink Dog"
$ws.Range("B169").Value = "1.0.0"
$ws.Range("C169").Value = 0.125
$ws.Range("D169").Value = "Fey"
$ws.Range("E169").Value = "Extraplanar, Forest, Grassland, Hills"
$ws.Range("F169").Value = "Magic Dogs"
$ws.Range("G169").Value = "None"
$ws.Range("H169").Value = "None"
$ws.Range("I169").Value = "None"

# Row heights matching the new rows' wrapped content.
$ws.Rows.Item(166).RowHeight = 75
$ws.Rows.Item(167).RowHeight = 75
$ws.Rows.Item(168).RowHeight = 90
$ws.Rows.Item(169).RowHeight = 30

# Dimension / view: user ended this session on the Creatures sheet, scrolled
# down to the newly-added rows.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 123
$ws.Range("D172").Select()

Write-Host "done"
